$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D-column price cells to remain text (avoid numeric auto-coercion)
# by temporarily applying a text number format, then clearing it after writing.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "66.538.23"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "3.248.89"
$ws.Range("E3").Value = "  +2.71%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "607.68"
$ws.Range("E5").Value = "  +0.88%  "

$ws.Range("D6").Value = "157.18"
$ws.Range("E6").Value = "  +2.13%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "3.249.82"
$ws.Range("E8").Value = "  +2.73%  "

$ws.Range("D9").Value = "0.551"
$ws.Range("E9").Value = "  +0.53%  "

$ws.Range("E10").Value = "  +2.68%  "

$ws.Range("D11").Value = "5.85"
$ws.Range("E11").Value = "  +6.28%  "

$ws.Range("D12").Value = "0.501"
$ws.Range("E12").Value = "  -2.99%  "

$ws.Range("E13").Value = "  +1.70%  "

$ws.Range("D14").Value = "39.14"
$ws.Range("E14").Value = "  +2.14%  "

$ws.Range("D15").Value = "3.783.71"
$ws.Range("E15").Value = "  +2.78%  "

$ws.Range("D16").Value = "66.649.97"
$ws.Range("E16").Value = "  +0.67%  "

$ws.Range("D17").Value = "7.44"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "3.249.82"
$ws.Range("E18").Value = "  +2.94%  "

$ws.Range("D19").Value = "0.114"
$ws.Range("E19").Value = "  +1.13%  "

$ws.Range("D20").Value = "507.23"
$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("D21").Value = "15.44"
$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").Value = "0.751"
$ws.Range("E22").Value = "  +3.38%  "

$ws.Range("D23").Value = "8.11"
$ws.Range("E23").Value = "  +0.74%  "

$ws.Range("D24").Value = "14.73"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").Value = "87.13"
$ws.Range("E25").Value = "  +3.23%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("E27").Value = "  +1.50%  "

$ws.Range("D28").Value = "9.08"
$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("D29").Value = "2.42"
$ws.Range("E29").Value = "  +1.02%  "

$ws.Range("D30").Value = "0.134"
$ws.Range("E30").Value = "  +51.83%  "

$ws.Range("D31").Value = "2.90"
$ws.Range("E31").Value = "  -5.13%  "

$ws.Range("D32").Value = "6.90"
$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("D33").Value = "28.06"
$ws.Range("E33").Value = "  +0.33%  "

$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("D35").Value = "1.15"
$ws.Range("E35").Value = "  -3.42%  "

$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").Value = "3.38"
$ws.Range("E37").Value = "  +21.59%  "

$ws.Range("D38").Value = "55.64"
$ws.Range("E38").Value = "  +1.74%  "

$ws.Range("D39").Value = "0.0₃0781"
$ws.Range("E39").Value = "  +15.86%  "

$ws.Range("D40").Value = "493.55"
$ws.Range("E40").Value = "  -1.40%  "

$ws.Range("E41").Value = "  +1.95%  "

$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("D43").Value = "8.83"
$ws.Range("E43").Value = "  +0.95%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.293"
$ws.Range("E44").Value = "  -0.61%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  +3.62%  "

$ws.Range("D46").Value = "2.983.71"
$ws.Range("E46").Value = "  +5.59%  "

$ws.Range("D47").Value = "28.86"
$ws.Range("E47").Value = "  +3.85%  "

$ws.Range("E48").Value = "  +6.37%  "

$ws.Range("E49").Value = "  +2.75%  "

$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").Value = "121.05"
$ws.Range("E51").Value = "  -0.45%  "

$dRange.ClearFormats()
